# Update FFXIV leve-profit calculation figures (currentAveragePrice / LevePrice / LeveProfit)
# for the rows whose underlying market-board prices changed, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5: Met a Sticky End / Animal Glue
$ws.Range("H5").Value = 3433.3333
$ws.Range("I5").Value = 5050
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 5050
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -4935
$ws.Range("N5").Value = -430

# Row 32: Automata for the People / Crab Oil
$ws.Range("H32").Value = 83336060
$ws.Range("I32").Value = 500000000
$ws.Range("J32").Value = 3280.4
$ws.Range("K32").Value = 500000000
$ws.Range("L32").Value = 3280.4
$ws.Range("M32").Value = -499999674
$ws.Range("N32").Value = -3932.4

$ws = $wb.Worksheets.Item("ARM")
# Row 33: A Leg to Stand On / Heavy Iron Flanchard
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 3594.2856
$ws.Range("I63").Value = 2678
$ws.Range("J63").Value = 4427.273
$ws.Range("K63").Value = 2678
$ws.Range("L63").Value = 4427.273
$ws.Range("M63").Value = -1992
$ws.Range("N63").Value = -5799.273

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 3594.2856
$ws.Range("I66").Value = 2678
$ws.Range("J66").Value = 4427.273
$ws.Range("K66").Value = 13390
$ws.Range("L66").Value = 22136.365
$ws.Range("M66").Value = -9958
$ws.Range("N66").Value = -29000.365

# Row 98: Greaving / Doman Iron Greaves of Maiming
$ws.Range("H98").Value = 38234
$ws.Range("J98").Value = 38234
$ws.Range("L98").Value = 38234
$ws.Range("N98").Value = -44224

# Row 112: Wrapped Knuckles / Deepgold Gloves of Fending
$ws.Range("H112").Value = 200031070
$ws.Range("J112").Value = 200031070
$ws.Range("L112").Value = 200031070
$ws.Range("N112").Value = -200034024

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 11112584
$ws.Range("I132").Value = 16130045
$ws.Range("J132").Value = 2492.7856
$ws.Range("K132").Value = 48390135
$ws.Range("L132").Value = 7478.3568
$ws.Range("M132").Value = -48387605
$ws.Range("N132").Value = -12538.3568

$ws = $wb.Worksheets.Item("BSM")
# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1425.6774
$ws.Range("I107").Value = 1306.1
$ws.Range("J107").Value = 5013
$ws.Range("K107").Value = 1306.1
$ws.Range("L107").Value = 5013
$ws.Range("M107").Value = 613.9000000000001
$ws.Range("N107").Value = -8853

# Row 122: To Delight a Dancer / High Durium Tathlums
$ws.Range("H122").Value = 40772
$ws.Range("J122").Value = 40772
$ws.Range("L122").Value = 40772
$ws.Range("N122").Value = -50572

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 238222.06
$ws.Range("I31").Value = 1807.3889
$ws.Range("J31").Value = 339542.66
$ws.Range("K31").Value = 1807.3889
$ws.Range("L31").Value = 339542.66
$ws.Range("M31").Value = -1512.3889
$ws.Range("N31").Value = -340132.66

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 238222.06
$ws.Range("I34").Value = 1807.3889
$ws.Range("J34").Value = 339542.66
$ws.Range("K34").Value = 1807.3889
$ws.Range("L34").Value = 339542.66
$ws.Range("M34").Value = -1605.3889
$ws.Range("N34").Value = -339946.66

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 1955.1666
$ws.Range("I99").Value = 1937.3334
$ws.Range("J99").Value = 1961.1111
$ws.Range("K99").Value = 1937.3334
$ws.Range("L99").Value = 1961.1111
$ws.Range("M99").Value = -439.3334
$ws.Range("N99").Value = -4957.1111

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 1955.1666
$ws.Range("I126").Value = 1937.3334
$ws.Range("J126").Value = 1961.1111
$ws.Range("K126").Value = 5812.0002
$ws.Range("L126").Value = 5883.3333
$ws.Range("M126").Value = -3342.0002
$ws.Range("N126").Value = -10823.3333

$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch / Chamomile Tea
$ws.Range("H34").Value = 1906.9286
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1906.9286
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5720.7858
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -5888.7858

# Row 35: Whirled Peas / Pea Soup
$ws.Range("H35").Value = 307.04544
$ws.Range("I35").Value = 1088.5
$ws.Range("J35").Value = 269.83334
$ws.Range("K35").Value = 3265.5
$ws.Range("L35").Value = 809.5000200000001
$ws.Range("M35").Value = -2977.5
$ws.Range("N35").Value = -1385.50002

# Row 36: Love's Crumpets Lost / Crumpet
$ws.Range("H36").Value = 300
$ws.Range("I36").Value = 300
$ws.Range("K36").Value = 900
$ws.Range("M36").Value = -731

# Row 39: Bloody Good Tart, This / Blood Currant Tart
$ws.Range("H39").Value = 1775.1666
$ws.Range("J39").Value = 1808.8695
$ws.Range("L39").Value = 5426.6085
$ws.Range("N39").Value = -6014.6085

# Row 55: Pagan Pastries / Pastry Fish
$ws.Range("H55").Value = 1234.6666
$ws.Range("I55").Value = 704
$ws.Range("J55").Value = 1500
$ws.Range("K55").Value = 2112
$ws.Range("L55").Value = 4500
$ws.Range("M55").Value = -1935
$ws.Range("N55").Value = -4854

# Row 129: Comfort Food / Yakow Moussaka
$ws.Range("H129").Value = 151768.6
$ws.Range("I129").Value = 429388.56
$ws.Range("J129").Value = 2280.923
$ws.Range("K129").Value = 1288165.68
$ws.Range("L129").Value = 6842.768999999999
$ws.Range("M129").Value = -1283165.68
$ws.Range("N129").Value = -16842.769

$ws = $wb.Worksheets.Item("GSM")
# Row 100: Hair-raising Action / Durium Hairpin of Maiming
$ws.Range("H100").Value = 36019
$ws.Range("J100").Value = 36019
$ws.Range("L100").Value = 36019
$ws.Range("N100").Value = -38183

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 2141.125
$ws.Range("I102").Value = 1928.7273
$ws.Range("J102").Value = 2608.4
$ws.Range("K102").Value = 1928.7273
$ws.Range("L102").Value = 2608.4
$ws.Range("M102").Value = -306.7273
$ws.Range("N102").Value = -5852.4

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 2287.3057
$ws.Range("I7").Value = 2041.8667
$ws.Range("J7").Value = 2462.6191
$ws.Range("K7").Value = 2041.8667
$ws.Range("L7").Value = 2462.6191
$ws.Range("M7").Value = -1929.8667
$ws.Range("N7").Value = -2686.6191

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 3900.4443
$ws.Range("J46").Value = 5250.6665
$ws.Range("L46").Value = 5250.6665
$ws.Range("N46").Value = -5626.6665

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 2770.8
$ws.Range("I61").Value = 2770.8
$ws.Range("K61").Value = 2770.8
$ws.Range("M61").Value = -2568.8

# Row 110: Breeches of Trust / Gliderskin Breeches of Fending
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 2770.8
$ws.Range("I113").Value = 2770.8
$ws.Range("K113").Value = 2770.8
$ws.Range("M113").Value = -600.8000000000002

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 2287.3057
$ws.Range("I126").Value = 2041.8667
$ws.Range("J126").Value = 2462.6191
$ws.Range("K126").Value = 6125.6001
$ws.Range("L126").Value = 7387.8573
$ws.Range("M126").Value = -3655.6001
$ws.Range("N126").Value = -12327.8573

# Row 133: The Perfect Accessory / Loboskin Amulet of Fending
$ws.Range("H133").Value = 34759.875
$ws.Range("J133").Value = 34759.875
$ws.Range("L133").Value = 34759.875
$ws.Range("N133").Value = -39819.875

Write-Output "Masamune_Profits: updated pricing rows across ALC/ARM/BSM/CRP/CUL/GSM/LTW"